$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Edad" (Age) column in F with header + values
$ws.Range("F1").Value = "Edad"
$ws.Range("F2").Value = 45
$ws.Range("F3").Value = 32
$ws.Range("F4").Value = 18
$ws.Range("F5").Value = 78

# Apply underline style to G4 (empty, styled cell) and select it
$ws.Range("G4").Font.Underline = $true
[void]$ws.Range("G4").Select()
